$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# ---------------------------------------------------------------------------
# 1. Build two reusable "template" cell styles off to the side, then copy
#    their formatting onto the real cells. This keeps the styles.xml /
#    borders list small (one style per unique combination) instead of
#    minting a new cellXf for every single cell we touch.
# ---------------------------------------------------------------------------

# "mid" template: thin top+bottom border, unlocked (protection)
$ws.Range("Z1").Value = "tmp"
$ws.Range("Z1").Locked = $false
$ws.Range("Z1").Borders.Item(8).Color = 0      # xlEdgeTop
$ws.Range("Z1").Borders.Item(8).LineStyle = 1
$ws.Range("Z1").Borders.Item(9).Color = 0      # xlEdgeBottom
$ws.Range("Z1").Borders.Item(9).LineStyle = 1

# "last column" template: thin top+bottom+right border, unlocked
$ws.Range("Z2").Value = "tmp2"
$ws.Range("Z2").Locked = $false
$ws.Range("Z2").Borders.Item(8).Color = 0      # xlEdgeTop
$ws.Range("Z2").Borders.Item(8).LineStyle = 1
$ws.Range("Z2").Borders.Item(9).Color = 0      # xlEdgeBottom
$ws.Range("Z2").Borders.Item(9).LineStyle = 1
$ws.Range("Z2").Borders.Item(10).Color = 0     # xlEdgeRight
$ws.Range("Z2").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------------
# 2. Header row (row 1) - plain text, default style (already default on this
#    template sheet: customFormat row, no per-cell style needed).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Order_Message"
$ws.Range("B1").Value = "Order No"
$ws.Range("C1").Value = "FlyFrom"
$ws.Range("D1").Value = "FlyTo"
$ws.Range("E1").Value = "Date"
$ws.Range("F1").Value = "Class"
$ws.Range("G1").Value = "Tickets"
$ws.Range("H1").Value = "Passenger"

# ---------------------------------------------------------------------------
# 3. Data rows 2-4.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Order 141 completed"
$ws.Range("C2").Value = "Denver"
$ws.Range("D2").Value = "Frankfurt"
$ws.Range("E2").Value = "'17-Jan-2021"
$ws.Range("F2").Value = "Business"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = "Amir khan"

$ws.Range("A3").Value = "Order 141 completed"
$ws.Range("C3").Value = "London"
$ws.Range("D3").Value = "Paris"
$ws.Range("E3").Value = "'15-Jan-2020"
$ws.Range("F3").Value = "Business"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "Roger Wattenhofer"

$ws.Range("A4").Value = "Order 141 completed"
$ws.Range("C4").Value = "Sydney"
$ws.Range("D4").Value = "Paris"
$ws.Range("E4").Value = "'17-Jan-2023"
$ws.Range("F4").Value = "Business"
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = "Samuel Mishler"

# Shared formula for the order number, across the whole B2:B4 block at once
# (mirrors Excel's own si="0" shared-formula grouping).
$ws.Range("B2:B4").Formula = "=SUBSTITUTE(SUBSTITUTE(A2, ""Order "", """"), "" completed"", """")"

# ---------------------------------------------------------------------------
# 4. Apply the border/protection formatting templates built in step 1.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Copy()
$ws.Range("A2:G4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("Z2").Copy()
$ws.Range("H2:H4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Clean up the scratch template cells.
# ---------------------------------------------------------------------------
$ws.Range("Z1:Z2").Clear()

# ---------------------------------------------------------------------------
# 6. Column widths (best effort / visual match).
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 18.36
$ws.Columns("B").ColumnWidth = 8.76
$ws.Columns("C").ColumnWidth = 9.42
$ws.Columns("D").ColumnWidth = 9.42
$ws.Columns("E").ColumnWidth = 10.98
$ws.Columns("F").ColumnWidth = 8.06
$ws.Columns("G").ColumnWidth = 6.86
$ws.Columns("H").ColumnWidth = 16.71

# ---------------------------------------------------------------------------
# 7. Selection, matching the recorded edit position.
# ---------------------------------------------------------------------------
$ws.Range("B3").Select()
